$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from serial 45178 (2023-09-09) to 45179 (2023-09-10)
# for all data rows (rows 2 through 33).
for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
